$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.04085365438401
$ws.Range("D2").Value = 1.052257361413211
$ws.Range("E2").Value = 1.049478834493663
$ws.Range("F2").Value = 1.061082028683717
$ws.Range("I2").Value = 1.042982717299858
$ws.Range("J2").Value = 1.045938016317031
$ws.Range("K2").Value = 1.055006594089548
$ws.Range("L2").Value = 1.05223577948633
$ws.Range("M2").Value = 1.063807056908867
$ws.Range("N2").Value = 1.019111224970098
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.042249780945378
$ws.Range("D3").Value = 1.052995776197247
$ws.Range("E3").Value = 1.050640395774867
$ws.Range("F3").Value = 1.062141771948873
$ws.Range("I3").Value = 1.043192547166726
$ws.Range("J3").Value = 1.046977704605015
$ws.Range("K3").Value = 1.055558009583848
$ws.Range("L3").Value = 1.05320869321399
$ws.Range("M3").Value = 1.064680733287245
$ws.Range("N3").Value = 1.019464255734456
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.043152564587135
$ws.Range("D4").Value = 1.053472073759189
$ws.Range("E4").Value = 1.05139147124112
$ws.Range("F4").Value = 1.062826302805115
$ws.Range("I4").Value = 1.04332594188614
$ws.Range("J4").Value = 1.047649423913213
$ws.Range("K4").Value = 1.055912641810683
$ws.Range("L4").Value = 1.05383712565092
$ws.Range("M4").Value = 1.065244267322451
$ws.Range("N4").Value = 1.019692125809572
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.043531954647306
$ws.Range("D5").Value = 1.053671948694716
$ws.Range("E5").Value = 1.051707098329339
$ws.Range("F5").Value = 1.063113795223161
$ws.Range("I5").Value = 1.043381451401179
$ws.Range("J5").Value = 1.047931571224356
$ws.Range("K5").Value = 1.056061209797688
$ws.Range("L5").Value = 1.054101055632682
$ws.Range("M5").Value = 1.065480749120798
$ws.Range("N5").Value = 1.019787788122297
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.043595647812622
$ws.Range("D6").Value = 1.05370548743265
$ws.Range("E6").Value = 1.051760086243647
$ws.Range("F6").Value = 1.063162049832708
$ws.Range("I6").Value = 1.043390738288698
$ws.Range("J6").Value = 1.047978930800788
$ws.Range("K6").Value = 1.056086124564237
$ws.Range("L6").Value = 1.054145355278103
$ws.Range("M6").Value = 1.065520430400881
$ws.Range("N6").Value = 1.019803842396955
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.043157634561162
$ws.Range("D7").Value = 1.053474745914661
$ws.Range("E7").Value = 1.051395689154076
$ws.Range("F7").Value = 1.062830145408432
$ws.Range("I7").Value = 1.043326685845319
$ws.Range("J7").Value = 1.047653194932746
$ws.Range("K7").Value = 1.055914629024094
$ws.Range("L7").Value = 1.053840653327035
$ws.Range("M7").Value = 1.065247428883613
$ws.Range("N7").Value = 1.019693404580594
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.041325609477219
$ws.Range("D8").Value = 1.052507224033498
$ws.Range("E8").Value = 1.04987150125024
$ws.Range("F8").Value = 1.061440421978086
$ws.Range("I8").Value = 1.043054122908556
$ws.Range("J8").Value = 1.046289598466005
$ws.Range("K8").Value = 1.055193397084167
$ws.Range("L8").Value = 1.052564810940426
$ws.Range("M8").Value = 1.064102691997766
$ws.Range("N8").Value = 1.019230650759599
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03809253018405
$ws.Range("D9").Value = 1.050790788308924
$ws.Range("E9").Value = 1.047181501254747
$ws.Range("F9").Value = 1.0589823404767
$ws.Range("I9").Value = 1.042555610284823
$ws.Range("J9").Value = 1.043878765984041
$ws.Range("K9").Value = 1.053905859712889
$ws.Range("L9").Value = 1.050308029840047
$ws.Range("M9").Value = 1.062071725539127
$ws.Range("N9").Value = 1.018410857631866
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.035933631459746
$ws.Range("D10").Value = 1.04963873069465
$ws.Range("E10").Value = 1.04538519772742
$ws.Range("F10").Value = 1.057337327552643
$ws.Range("I10").Value = 1.042211014491157
$ws.Range("J10").Value = 1.042265988777262
$ws.Range("K10").Value = 1.053036297116612
$ws.Range("L10").Value = 1.048797595849952
$ws.Range("M10").Value = 1.060708383424948
$ws.Range("N10").Value = 1.017861340865055
$ws.Range("B11").Value = 1.019999999999999
$ws.Range("C11").Value = 1.03499790971204
$ws.Range("D11").Value = 1.049138031175774
$ws.Range("E11").Value = 1.044606638344901
$ws.Range("F11").Value = 1.056623503074934
$ws.Range("I11").Value = 1.042058892712587
$ws.Range("J11").Value = 1.04156628077834
$ws.Range("K11").Value = 1.052657103557216
$ws.Range("L11").Value = 1.048142126185659
$ws.Range("M11").Value = 1.060115798724604
$ws.Range("N11").Value = 1.01762267266412
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.034650199720728
$ws.Range("D12").Value = 1.048951770522201
$ws.Range("E12").Value = 1.044317331148598
$ws.Range("F12").Value = 1.056358125899953
$ws.Range("I12").Value = 1.042001950547419
$ws.Range("J12").Value = 1.04130616948025
$ws.Range("K12").Value = 1.052515852814261
$ws.Range("L12").Value = 1.047898436043748
$ws.Range("M12").Value = 1.059895346570144
$ws.Range("N12").Value = 1.017533910790249
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034724791174095
$ws.Range("D13").Value = 1.048991736677905
$ws.Range("E13").Value = 1.044379393805741
$ws.Range("F13").Value = 1.056415060692268
$ws.Range("I13").Value = 1.042014184630785
$ws.Range("J13").Value = 1.041361973733087
$ws.Range("K13").Value = 1.052546169751762
$ws.Range("L13").Value = 1.047950718387367
$ws.Range("M13").Value = 1.059942649715084
$ws.Range("N13").Value = 1.01755295550396
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034969170808942
$ws.Range("D14").Value = 1.04912264049739
$ws.Range("E14").Value = 1.044582726498034
$ws.Range("F14").Value = 1.056601571639433
$ws.Range("I14").Value = 1.042054194778981
$ws.Range("J14").Value = 1.041544784173116
$ws.Range("K14").Value = 1.052645435918465
$ws.Range("L14").Value = 1.048121987177669
$ws.Range("M14").Value = 1.060097583024424
$ws.Range("N14").Value = 1.017615337827982
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035119722285557
$ws.Range("D15").Value = 1.04920325773462
$ws.Range("E15").Value = 1.044707991069929
$ws.Range("F15").Value = 1.056716456517151
$ws.Range("I15").Value = 1.042078788389424
$ws.Range("J15").Value = 1.041657391971718
$ws.Range("K15").Value = 1.052706543874391
$ws.Range("L15").Value = 1.048227482355349
$ws.Range("M15").Value = 1.060192997461338
$ws.Range("N15").Value = 1.017653759042852
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.035995712716339
$ws.Range("D16").Value = 1.049671921435336
$ws.Range("E16").Value = 1.045436852116861
$ws.Range("F16").Value = 1.057384669459942
$ws.Range("I16").Value = 1.042221048992131
$ws.Range("J16").Value = 1.042312397035574
$ws.Range("K16").Value = 1.053061406694048
$ws.Range("L16").Value = 1.048841066588168
$ws.Range("M16").Value = 1.060747663775952
$ws.Range("N16").Value = 1.017877165141119
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.036544952444678
$ws.Range("D17").Value = 1.049965405911711
$ws.Range("E17").Value = 1.045893844868679
$ws.Range("F17").Value = 1.057803412509545
$ws.Range("I17").Value = 1.042309506182661
$ws.Range("J17").Value = 1.042722896510843
$ws.Range("K17").Value = 1.05328328822331
$ws.Range("L17").Value = 1.049225563712749
$ws.Range("M17").Value = 1.061094987899417
$ws.Range("N17").Value = 1.018017107301657
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.03686522779634
$ws.Range("D18").Value = 1.050136411884502
$ws.Range("E18").Value = 1.046160329045726
$ws.Range("F18").Value = 1.058047511445047
$ws.Range("I18").Value = 1.042360821056114
$ws.Range("J18").Value = 1.042962202465412
$ws.Range("K18").Value = 1.053412450664278
$ws.Range("L18").Value = 1.04944969546867
$ws.Range("M18").Value = 1.061297359367388
$ws.Range("N18").Value = 1.018098663435963
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.03697441883918
$ws.Range("D19").Value = 1.050194690228985
$ws.Range("E19").Value = 1.0462511811527
$ws.Range("F19").Value = 1.058130718038459
$ws.Range("I19").Value = 1.042378270495448
$ws.Range("J19").Value = 1.043043777422885
$ws.Range("K19").Value = 1.053456448094953
$ws.Range("L19").Value = 1.049526095130849
$ws.Range("M19").Value = 1.061366326091567
$ws.Range("N19").Value = 1.018126460171582
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.036486033247826
$ws.Range("D20").Value = 1.049933936278051
$ws.Range("E20").Value = 1.045844821339098
$ws.Range("F20").Value = 1.05775850054433
$ws.Range("I20").Value = 1.042300044600464
$ws.Range("J20").Value = 1.042678867424799
$ws.Range("K20").Value = 1.053259509054346
$ws.Range("L20").Value = 1.049184325200586
$ws.Range("M20").Value = 1.061057745749821
$ws.Range("N20").Value = 1.018002100056473
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.034897211040803
$ws.Range("D21").Value = 1.049084100278754
$ws.Range("E21").Value = 1.044522853304253
$ws.Range("F21").Value = 1.056546655207888
$ws.Range("I21").Value = 1.042042424862945
$ws.Range("J21").Value = 1.041490956841724
$ws.Range("K21").Value = 1.052616215591699
$ws.Range("L21").Value = 1.048071558871138
$ws.Range("M21").Value = 1.060051968423034
$ws.Range("N21").Value = 1.017596970830686
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033897435409231
$ws.Range("D22").Value = 1.048548162747359
$ws.Range("E22").Value = 1.043691009247546
$ws.Range("F22").Value = 1.055783382468681
$ws.Range("I22").Value = 1.041877918528418
$ws.Range("J22").Value = 1.040742861683277
$ws.Range("K22").Value = 1.052209428906149
$ws.Range("L22").Value = 1.047370646956312
$ws.Range("M22").Value = 1.059417629244778
$ws.Range("N22").Value = 1.017341613723446
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.034427514930324
$ws.Range("D23").Value = 1.048832426327756
$ws.Range("E23").Value = 1.044132050125198
$ws.Range("F23").Value = 1.056188135296511
$ws.Range("I23").Value = 1.041965366402939
$ws.Range("J23").Value = 1.041139556834105
$ws.Range("K23").Value = 1.052425294600236
$ws.Range("L23").Value = 1.047742335131581
$ws.Range("M23").Value = 1.059754091570676
$ws.Range("N23").Value = 1.017477044049925
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.036512656564893
$ws.Range("D24").Value = 1.04994815660275
$ws.Range("E24").Value = 1.045866973184173
$ws.Range("F24").Value = 1.057778794781787
$ws.Range("I24").Value = 1.042304320749589
$ws.Range("J24").Value = 1.042698762679542
$ws.Range("K24").Value = 1.053270254632332
$ws.Range("L24").Value = 1.049202959538244
$ws.Range("M24").Value = 1.061074574544065
$ws.Range("N24").Value = 1.018008881400097
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.038928957615875
$ws.Range("D25").Value = 1.051235896306691
$ws.Range("E25").Value = 1.047877442749789
$ws.Range("F25").Value = 1.059618915106541
$ws.Range("I25").Value = 1.042686647345533
$ws.Range("J25").Value = 1.04450299019152
$ws.Range("K25").Value = 1.054240692245014
$ws.Range("L25").Value = 1.050892493607058
$ws.Range("M25").Value = 1.062598422878149
$ws.Range("N25").Value = 1.018623316373808
